# fix: coba tawarkan ke ustadz wahidi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 10 (16 Maret) as confirmed
$ws.Range("E10").Value = "✔️"

# Offer the slot to Ustadz Wahidi on 18 Maret (previously Ustadz Miftah)
$ws.Range("D12").Value = "Ustadz Wahidi"

# Offer the open slot on 21 Maret to Ustadz Wahidi as well
$ws.Range("D15").Value = "Ustadz Wahidi"

# Update the view's selection to D16 and clear the scrolled topLeftCell
$ws.Range("D16").Select()
